# Daily attendance processing - reorder "Recorded By" (column G) entries
# so that any "System"/"system" token is listed first, by reversing the
# comma-separated list of recorders for every data row that contains one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$processed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -eq $null) { continue }
    if ($v -notlike "*,*") { continue }
    if ($v -notlike "*system*") { continue }

    $parts = $v -split ", "
    $reversed = $parts[($parts.Length - 1)..0]
    $joined = [string]::Join(", ", $reversed)

    $cell.Value = $joined
    $processed++
}

Write-Host ("Reordered Recorded By for " + $processed + " rows")
